# Atualizacao de bases das ligas: alguns jogos tinham sido gravados nas
# linhas trocadas - troca-se o conteudo (colunas B:AC) entre as linhas
# pares correspondentes, mantendo o indice da coluna A inalterado.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estonia Meistriliiga")

function Swap-MatchRows {
    param($ws, [int]$row1, [int]$row2, [int]$firstCol, [int]$lastCol)

    # Only write back cells whose value actually differs between the
    # two rows (columns already equal on both rows are left untouched).
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($v1 -ne $v2) {
            $cell1.Value2 = $v2
            $cell2.Value2 = $v1
        }
    }
}

Swap-MatchRows $ws 71 72 2 29
Swap-MatchRows $ws 104 107 2 29
Swap-MatchRows $ws 105 106 2 29
